# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Swap the "Periodo Mora" values for the two worker rows (E16/E17):
#     E16: 1801 -> 1703
#     E17: 1703 -> 1801
# - Update "Salario Basico" (column G) for both rows with the new amount:
#     G16: 737717 -> 781242
#     G17: 737717 -> 781242

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Periodo Mora (text-formatted period codes) - values swap between the two rows
$ws.Range("E16").Value = "1703"
$ws.Range("E17").Value = "1801"

# Salario Basico - updated amount for both rows
$ws.Range("G16").Value = 781242
$ws.Range("G17").Value = 781242
